$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.137.20"
$ws.Range("D3").Value = "2.047.32"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.08"
$ws.Range("E5").Value = "  -1.96%  "
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.94"
$ws.Range("E8").Value = "  -2.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.382"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0781"
$ws.Range("E10").Value = "  -2.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.109"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.20"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.896"
$ws.Range("E13").Value = "  +10.67%  "
$ws.Range("D14").Value = "2.345.71"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.71"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").Value = "2.048.63"
$ws.Range("E16").Value = "  -1.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.58"
$ws.Range("E17").Value = "  +13.07%  "
$ws.Range("D18").Value = "37.134.09"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.39"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("D20").Value = "0.0₃0893"
$ws.Range("E20").Value = "  -3.41%  "
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "236.46"
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.49"
$ws.Range("E24").Value = "  +3.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.58"
$ws.Range("E25").Value = "  +2.51%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.18"
$ws.Range("E26").Value = "  -4.45%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.67"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.15"
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("E29").Value = "  -1.23%  "
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.96"
$ws.Range("E31").Value = "  +3.64%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0885"
$ws.Range("E34").Value = "  -3.30%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.25"
$ws.Range("E36").Value = "  -1.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.79"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("E38").Value = "  -2.38%  "
$ws.Range("E39").Value = "  +13.00%  "
$ws.Range("E40").Value = "  +8.37%  "
$ws.Range("E41").Value = "  -15.35%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0223"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.53"
$ws.Range("E43").Value = "  -2.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "95.79"
$ws.Range("E45").Value = "  -2.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.43"
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("D47").Value = "1.270.55"
$ws.Range("E47").Value = "  -2.14%  "
$ws.Range("E48").Value = "  -2.38%  "
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("D50").Value = "2.230.72"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("E51").Value = "  +0.34%  "
